# Increment the "Förändrad" (changed) date in column C for rows 2-9
# from 45183 (2023-09-14) to 45184 (2023-09-15), matching the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
